$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("period_lbl")

# Column D ("date_added"), rows 2 through 145, all get updated to the new
# upload date ("2025-01-15") for the new-format backdata upload.
#
# Assigning the literal "2025-01-15" string straight to .Value/.Formula gets
# auto-recognized as a date and silently converted to a date serial (plus a
# new number-format style) instead of staying plain text like the source
# data. To keep it as literal text (matching the existing date_added
# entries), write it as a text formula and then flatten the range down to
# plain values via copy / paste-special-values - this keeps the cells as
# literal strings with no extra number formatting applied.
$targetRange = $ws.Range("D2:D145")
$targetRange.Formula = '="2025-01-15"'
$targetRange.Copy()
$targetRange.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false
